# TestPlan.xlsx update:
# 1. New test cases added for Product Page (Start Today / Get Started)
# 2. Updated test case for "Become a Partner" page (now part of test plan)
# 3. "Debug tests" row is no longer part of the test plan

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---------------------------------------------
# Row 2 ("Become a Partner" / BecomePartnerPageTestCases.xlsx): N -> Y
$ws.Range("D2").Value = "Y"

# Row 20 ("Debug tests" / Test.xlsx): Y -> N
$ws.Range("D20").Value = "N"

# --- Append new rows for the Product Page test cases -------------------
# Column A (Sno)
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(22, 1).Value = 21

# Column B (Test Suite Name)
$ws.Cells.Item(21, 2).Value = "Product Page Start Today"
$ws.Cells.Item(22, 2).Value = "Product Page Get Started"

# Column C (Target Page) - copy formatting (Menlo style) from an existing cell
$ws.Cells.Item(19, 3).Copy()
$ws.Cells.Item(21, 3).PasteSpecial(-4122)
$ws.Cells.Item(19, 3).Copy()
$ws.Cells.Item(22, 3).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(22, 3).Value = "ProductsCapellaGetstartedPageTestCases.xlsx"
$ws.Cells.Item(21, 3).Value = "ProductsCapellaStartTodayPageTestCases.xlsx"

# Column D (Part Of Test Plan)
$ws.Cells.Item(21, 4).Value = "Y"
$ws.Cells.Item(22, 4).Value = "Y"

# --- Update selection to match the post-edit view (active cell D2) -----
$ws.Range("D2").Select() | Out-Null
